$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The old row 1 is a blank spacer row above the header; delete it so the
# header row moves up to row 1 and everything below shifts up by one.
$ws.Rows.Item(1).Delete()

# The body rows (now rows 3-11, columns B:E) carried a stray "apply border"
# flag left over from an earlier edit even though the border itself is
# invisible (borderId 0). Clear any border formatting there so those cells
# fall back to the plain right/vertical-center format already used
# elsewhere in the sheet.
$ws.Range("B3:E11").Borders.LineStyle = -4142

# Select C15 on the sheet (the cell below the now-shifted table), matching
# where the cursor ends up after the row delete.
$ws.Range("C15").Select()
